$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "comments / blue / red / purple" row (previously row 8) moves up to row 2,
# and the three file rows shift down by one (row 2->3, 3->4, 4->5).
$ws.Range("A2").Value = "comments"
$ws.Range("B2").Value = "blue"
$ws.Range("C2").Value = "red"
$ws.Range("D2").Value = "purple "

$ws.Range("A3").Value = "main.py"
$ws.Range("B3").Value = "https://drive.google.com/file/d/1C3uKfNAam_jvbCEko8mAGleKwRaLTiBH/view?usp=drive_link"
$ws.Range("C3").Value = "https://drive.google.com/file/d/1Ls6ffkKL4olnW9JitVJU19WV5jim_NLM/view?usp=drive_link"
$ws.Range("D3").Value = "https://drive.google.com/file/d/1FYspWxD-uI9dgEIrUKxZ_IVb_X24ssBo/view?usp=drive_link"

$ws.Range("A4").Value = "main_out.mpy"
$ws.Range("B4").Value = "https://drive.google.com/file/d/1Iho4UCAvmH5JCXC1XlUxdqt5aso4QvJf/view?usp=drive_link"
$ws.Range("C4").Value = "https://drive.google.com/file/d/1JloC9vFhQ92tlrtBBpJFwczeOmq_kJfT/view?usp=drive_link"
$ws.Range("D4").Value = "https://drive.google.com/file/d/1q_WCGAmE2CNaKRg6wlXEt2QX7qp-25yC/view?usp=drive_link"

$ws.Range("A5").Value = "pio_ws2812_obj.mpy"
$ws.Range("B5").Value = "https://drive.google.com/file/d/1XacaFMCTDehv7wjFZgKxgRFxwUHD9vv-/view?usp=drive_link"
$ws.Range("C5").Value = "https://drive.google.com/file/d/1os18rWyXREjMl9tCmJM40SwLB7DUlXwA/view?usp=drive_link"
$ws.Range("D5").Value = "https://drive.google.com/file/d/1bJvsO_phaL0UM1HKEAYpxogippoHLXpT/view?usp=drive_link"

# Remove the now-duplicate/old data from row 8 (it already moved up to row 2).
$ws.Range("A8:D8").ClearContents()

# New column widths for the updated (narrower) control.
$ws.Columns.Item(1).ColumnWidth = 18.428571428571427
$ws.Columns.Item(2).ColumnWidth = 15.571428571428571
$ws.Columns.Item(3).ColumnWidth = 13.285714285714286
$ws.Columns.Item(4).ColumnWidth = 11.285714285714286

# Update the active selection shown in the saved view.
$ws.Range("F9").Select()
